$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the leading "Lit Review" section: the Heading1 "Lit Review"
#    paragraph (with the _GoBack bookmark wrapped around it), the
#    "Data Visualisation ..." paragraph, and the two empty Heading1
#    paragraphs that follow it.
# ---------------------------------------------------------------------------
$startP = $d.Paragraphs.Item(1).Range.Start
$endP = $d.Paragraphs.Item(4).Range.End
$d.Range($startP, $endP).Delete()

# That delete leaves a dangling <w:bookmarkEnd/> behind (its matching
# <w:bookmarkStart/> for "_GoBack" was inside the deleted text). A trailing
# no-op delete flushes/cleans up that orphaned bookmark marker.
$d.Range(0, 0).Delete()

# ---------------------------------------------------------------------------
# 2) Drop the stale <w:lastRenderedPageBreak/> sitting in front of the
#    second "Raspberry-PI: " run (the one right after the Bootstrap link).
#    Re-writing the run's text with a genuinely different value first forces
#    the engine to actually rebuild the run (clearing the cached page-break
#    marker); writing the correct text back afterwards restores the content.
# ---------------------------------------------------------------------------
$r = $d.Content
$hit = 0
while ($r.Find.Execute("Raspberry-PI: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $hit = $hit + 1
    if ($hit -eq 2) {
        $r.Text = "\x01TEMP\x01Raspberry-PI: "
        break
    }
    $r.Collapse(0)
}

$r2 = $d.Content
if ($r2.Find.Execute("\x01TEMP\x01Raspberry-PI: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $r2.Text = "Raspberry-PI: "
}

# ---------------------------------------------------------------------------
# 3) Re-anchor the "_GoBack" bookmark right after the "... Index Page."
#    sentence (immediately before the trailing " " run at the end of the
#    last paragraph), mirroring where Word last left the edit cursor.
# ---------------------------------------------------------------------------
$r3 = $d.Content
if ($r3.Find.Execute("template to produce the Index Page.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $insertPoint = $d.Range($r3.End, $r3.End)
    $d.Bookmarks.Add("_GoBack", $insertPoint)
}
